# Auto-generated edit script: scheduled market-data refresh
# Updates literal price/profit columns (H:N) across multiple sheets
# to reflect refreshed Universalis market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 883
$ws.Range("I2").Value = 879.8
$ws.Range("K2").Value = 879.8
$ws.Range("M2").Value = -766.8
$ws.Range("H16").Value = 762.4
$ws.Range("I16").Value = 762.4
$ws.Range("K16").Value = 762.4
$ws.Range("M16").Value = -532.4
$ws.Range("H70").Value = 3145.5386
$ws.Range("I70").Value = 2710.2222
$ws.Range("J70").Value = 4125
$ws.Range("K70").Value = 8130.6666
$ws.Range("L70").Value = 12375
$ws.Range("M70").Value = -7860.6666
$ws.Range("N70").Value = -12915
$ws.Range("H73").Value = 3145.5386
$ws.Range("I73").Value = 2710.2222
$ws.Range("J73").Value = 4125
$ws.Range("K73").Value = 8130.6666
$ws.Range("L73").Value = 12375
$ws.Range("M73").Value = -7194.6666
$ws.Range("N73").Value = -14247
$ws.Range("H94").Value = 12502.5
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H98").Value = 1451.909
$ws.Range("I98").Value = 1451.909
$ws.Range("K98").Value = 1451.909
$ws.Range("M98").Value = 46.09099999999989
$ws.Range("H105").Value = 9951.25
$ws.Range("J105").Value = 9951.25
$ws.Range("L105").Value = 9951.25
$ws.Range("N105").Value = -16939.25
$ws.Range("H107").Value = 65028.645
$ws.Range("I107").Value = 65028.645
$ws.Range("K107").Value = 65028.645
$ws.Range("M107").Value = -63108.645
$ws.Range("H116").Value = 5002.5
$ws.Range("J116").Value = 5002.5
$ws.Range("L116").Value = 5002.5
$ws.Range("N116").Value = -11886.5
$ws.Range("H122").Value = 1451.909
$ws.Range("I122").Value = 1451.909
$ws.Range("K122").Value = 4355.727000000001
$ws.Range("M122").Value = -1905.727000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1050
$ws.Range("I5").Value = 900
$ws.Range("J5").Value = 1125
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 1125
$ws.Range("M5").Value = -788
$ws.Range("N5").Value = -1349
$ws.Range("H25").Value = 5986.375
$ws.Range("I25").Value = 3000
$ws.Range("J25").Value = 6981.8335
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 6981.8335
$ws.Range("M25").Value = -2598
$ws.Range("N25").Value = -7785.8335
$ws.Range("H31").Value = 1999
$ws.Range("I31").Value = 1999
$ws.Range("K31").Value = 1999
$ws.Range("M31").Value = -1705
$ws.Range("H102").Value = 23334612
$ws.Range("I102").Value = 1251438
$ws.Range("K102").Value = 1251438
$ws.Range("M102").Value = -1249816

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1050
$ws.Range("I4").Value = 900
$ws.Range("J4").Value = 1125
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1125
$ws.Range("M4").Value = -785
$ws.Range("N4").Value = -1355
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("H22").Value = 760.4545000000001
$ws.Range("I22").Value = 707.3333
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 707.3333
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -534.3333
$ws.Range("N22").Value = -1345.5
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -748
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("N42").Value = 0
$ws.Range("H82").Value = 33702.4
$ws.Range("I82").Value = 20378.125
$ws.Range("J82").Value = 86999.5
$ws.Range("K82").Value = 20378.125
$ws.Range("L82").Value = 86999.5
$ws.Range("M82").Value = -19995.125
$ws.Range("N82").Value = -87765.5
$ws.Range("H85").Value = 33702.4
$ws.Range("I85").Value = 20378.125
$ws.Range("J85").Value = 86999.5
$ws.Range("K85").Value = 20378.125
$ws.Range("L85").Value = 86999.5
$ws.Range("M85").Value = -19052.125
$ws.Range("N85").Value = -89651.5
$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -2502
$ws.Range("H135").Value = 130000
$ws.Range("J135").Value = 130000
$ws.Range("L135").Value = 130000
$ws.Range("N135").Value = -140140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 861332
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 1066665
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 1066665
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -1077025

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 32349.621
$ws.Range("I34").Value = 162.71428
$ws.Range("J34").Value = 42590.91
$ws.Range("K34").Value = 488.14284
$ws.Range("L34").Value = 127772.73
$ws.Range("M34").Value = -404.14284
$ws.Range("N34").Value = -127940.73
$ws.Range("H106").Value = 3425
$ws.Range("J106").Value = 3425
$ws.Range("L106").Value = 10275
$ws.Range("N106").Value = -12167
$ws.Range("H137").Value = 999.6667
$ws.Range("J137").Value = 999.6667
$ws.Range("L137").Value = 2999.0001
$ws.Range("N137").Value = -13199.0001
$ws.Range("H140").Value = 743.3333
$ws.Range("I140").Value = 692
$ws.Range("K140").Value = 2076
$ws.Range("M140").Value = 3104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 399.6154
$ws.Range("I2").Value = 489.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 489.5
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -376.5
$ws.Range("N2").Value = -326
$ws.Range("H9").Value = 335.66666
$ws.Range("I9").Value = 335.66666
$ws.Range("K9").Value = 335.66666
$ws.Range("M9").Value = -165.66666
$ws.Range("H43").Value = 21329.8
$ws.Range("I43").Value = 824.5
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 824.5
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -673.5
$ws.Range("N43").Value = -35302
$ws.Range("H46").Value = 40682
$ws.Range("J46").Value = 40682
$ws.Range("L46").Value = 40682
$ws.Range("N46").Value = -40994
$ws.Range("H57").Value = 50060
$ws.Range("J57").Value = 50060
$ws.Range("L57").Value = 50060
$ws.Range("N57").Value = -51700

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 126.666664
$ws.Range("I9").Value = 126.666664
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 126.666664
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 97.333336
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0
$ws.Range("H22").Value = 214.08333
$ws.Range("I22").Value = 189.875
$ws.Range("J22").Value = 262.5
$ws.Range("K22").Value = 189.875
$ws.Range("L22").Value = 262.5
$ws.Range("M22").Value = 105.125
$ws.Range("N22").Value = -852.5
$ws.Range("H27").Value = 214.08333
$ws.Range("I27").Value = 189.875
$ws.Range("J27").Value = 262.5
$ws.Range("K27").Value = 189.875
$ws.Range("L27").Value = 262.5
$ws.Range("M27").Value = -82.875
$ws.Range("N27").Value = -476.5
$ws.Range("H30").Value = 828.75
$ws.Range("I30").Value = 828.75
$ws.Range("K30").Value = 828.75
$ws.Range("M30").Value = -720.75
$ws.Range("H35").Value = 11051.667
$ws.Range("I35").Value = 1494
$ws.Range("J35").Value = 22998.75
$ws.Range("K35").Value = 1494
$ws.Range("L35").Value = 22998.75
$ws.Range("M35").Value = -1158
$ws.Range("N35").Value = -23670.75
$ws.Range("H93").Value = 33334822
$ws.Range("I93").Value = 37038530
$ws.Range("J93").Value = 1484
$ws.Range("K93").Value = 37038530
$ws.Range("L93").Value = 1484
$ws.Range("M93").Value = -37037282
$ws.Range("N93").Value = -3980
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0
